$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.490.40"
$ws.Range("E2").Value = "  +2.95%  "
$ws.Range("D3").Value = "'1.607.14"
$ws.Range("E3").Value = "  +2.81%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'212.81"
$ws.Range("E5").Value = "  +1.05%  "
$ws.Range("E6").Value = "  +6.92%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "'26.92"
$ws.Range("E8").Value = "  +7.07%  "
$ws.Range("D9").Value = "'43.56"
$ws.Range("E9").Value = "  -0.84%  "
$ws.Range("E10").Value = "  +2.70%  "
$ws.Range("E11").Value = "  +2.76%  "
$ws.Range("D12").Value = "'0.0910"
$ws.Range("E12").Value = "  +1.56%  "
$ws.Range("D13").Value = "'1.837.23"
$ws.Range("E13").Value = "  +2.86%  "
$ws.Range("D14").Value = "'1.614.95"
$ws.Range("E14").Value = "  +3.10%  "
$ws.Range("D15").Value = "'29.501.98"
$ws.Range("E15").Value = "  +3.05%  "
$ws.Range("E16").Value = "  +3.98%  "
$ws.Range("E17").Value = "  +2.02%  "
$ws.Range("D18").Value = "'63.30"
$ws.Range("E18").Value = "  +3.32%  "
$ws.Range("D19").Value = "'241.30"
$ws.Range("E19").Value = "  +5.13%  "
$ws.Range("E20").Value = "  +3.69%  "
$ws.Range("E21").Value = "  +2.03%  "
$ws.Range("D22").Value = "'0.999"
$ws.Range("E22").Value = "  -0.07%  "
$ws.Range("E23").Value = "  +2.64%  "
$ws.Range("E24").Value = "  +2.26%  "
$ws.Range("E25").Value = "  +0.39%  "
$ws.Range("D26").Value = "'154.54"
$ws.Range("E26").Value = "  +2.41%  "
$ws.Range("E27").Value = "  +5.02%  "
$ws.Range("E28").Value = "  +3.24%  "
$ws.Range("E29").Value = "  +2.48%  "
$ws.Range("D30").Value = "'1.00"
$ws.Range("E30").Value = "  -0.05%  "
$ws.Range("E31").Value = "  +2.58%  "
$ws.Range("D32").Value = "'1.07"
$ws.Range("E32").Value = "  +1.33%  "
$ws.Range("E33").Value = "  +1.61%  "
$ws.Range("E34").Value = "  +4.13%  "
$ws.Range("D35").Value = "'1.412.70"
$ws.Range("E35").Value = "  +1.83%  "
$ws.Range("E36").Value = "  +0.60%  "
$ws.Range("E37").Value = "  +3.50%  "
$ws.Range("D38").Value = "'2.82"
$ws.Range("E38").Value = "  +5.22%  "
$ws.Range("E39").Value = "  +0.25%  "
$ws.Range("E40").Value = "  +2.73%  "
$ws.Range("D41").Value = "'0.537"
$ws.Range("E41").Value = "  +3.75%  "
$ws.Range("E42").Value = "  +1.10%  "
$ws.Range("D43").Value = "'0.0485"
$ws.Range("E43").Value = "  +5.23%  "
$ws.Range("E44").Value = "  +3.48%  "
$ws.Range("E45").Value = "  -0.05%  "
$ws.Range("D46").Value = "'52.88"
$ws.Range("E46").Value = "  +22.24%  "
$ws.Range("D47").Value = "'65.67"
$ws.Range("E47").Value = "  +2.56%  "
$ws.Range("E48").Value = "  +0.94%  "
$ws.Range("D49").Value = "'1.748.53"
$ws.Range("E50").Value = "  -0.99%  "
$ws.Range("D51").Value = "'86.71"
$ws.Range("E51").Value = "  +1.92%  "
